# Correccion SRS: mark several review comments as "resolved" by striking
# them through (and, for one paragraph, swap a yellow highlight for strike).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Simple cases: apply strike-through to the whole paragraph (text +
#    paragraph mark) for paragraphs that currently have no strike at all.
# ---------------------------------------------------------------------
$simpleTargets = @(
    "Alcance",
    "Restricciones",
    "Suposiciones y dependencias"
)

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    foreach ($needle in $simpleTargets) {
        if ($text -like "*$needle*") {
            $p.Range.Font.StrikeThrough = 1
        }
    }
}

# ---------------------------------------------------------------------
# 2) "Interfaces de software": drop the yellow highlight and strike the
#    text instead.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    if ($text -like "*Interfaces de software*") {
        $p.Range.Font.StrikeThrough = 1
        $p.Range.HighlightColorIndex = 0
    }
}

# ---------------------------------------------------------------------
# 3) "Definiciones, acronimos y abreviaturas": the paragraph is split in
#    two runs ("...completar " without strike, "y ordenar
#    alfabeticamente" with strike). Strike the whole paragraph, then
#    merge the two runs into a single one, since the word "ordenar" no
#    longer needs its own distinct run now that both share formatting.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    if ($text -like "*Definiciones*") {
        $p.Range.Font.StrikeThrough = 1

        $pStart = $p.Range.Start
        $pEnd = $p.Range.End  # includes the trailing paragraph mark

        $searchRange = $d.Range($pStart, $pEnd)
        $searchRange.Find.ClearFormatting()
        $searchRange.Find.Execute("completar ")
        $boundary = $searchRange.End

        $tailRange = $d.Range($boundary, $pEnd - 1)
        $tailText = $tailRange.Text
        $tailRange.Delete()

        $insertPoint = $d.Range($boundary, $boundary)
        $insertPoint.InsertAfter($tailText)
    }
}
